# Apply the table-style change on slide 16's table (graphicFrame / Shapes.Item(3))
# and re-point the presentation's live theme colour scheme to the "Office Theme"
# palette (the deck's two theme parts effectively swap roles in the target edit;
# the single colour scheme this runtime exposes is repainted with the Office
# Theme's 12 colours, matching the values that the target theme XML ends up with).

$p = $ppt.ActivePresentation

# --- 1. Table style id swap on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{6E73156F-5722-4AA0-9FE4-88CC6D680518}")
}

# --- 2. Theme colour scheme -> Office Theme palette -------------------------
function ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = ComRGB($officeThemeColors[$i])
}
